# Add team record (Wins / Losses / Ties) columns to the WSN_2019 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - reuse the bold/bordered/centered header style (s="1")
# that the rest of row 1 already uses, by copying the format from an
# existing header cell.
$headerStyleCell = $ws.Range("AC1")

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerStyleCell.Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-52: constant team record for every player.
$lastRow = 52
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 93   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 69   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
